$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commercial Invoice")

# Clear the "Item" (C) and "Description" (D) columns for rows 9-14 to an
# empty string, and zero out Qty (E), Unit Price (F) and Unit Discount (G)
# for the same rows, matching the target state of the InvoiceItems table.
for ($r = 9; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
}
